$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename column H header from "Notebook Upload Path" to "Payment Confirmed"
$ws.Range("H1").Value = "Payment Confirmed"

# Add a new (blank, but formatted) cell in I1, matching the header style
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("I1").Value = $null

# Move the active selection (cosmetic, matches author's last cursor position)
$ws.Range("J21").Select()
